# "Get Test Data Via Util Library"
#
# Rename the "Fname" header to "FirstName" and drop the yellow header
# highlight from the B1:D1 cells (A1 keeps its highlight), then restore
# the view state (active cell / column A width) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the second header cell's text: "Fname" -> "FirstName"
$ws.Range("B1").Value = "FirstName"

# Remove the yellow fill from the B1:D1 header cells (A1's highlight is
# left untouched).
$ws.Range("B1:D1").Interior.Pattern = -4142

# Restore cosmetic view state: active cell moves to D1 ...
$null = $ws.Range("D1").Select()

# ... and column A widens slightly.
$ws.Columns("A").ColumnWidth = 17.166666666666668
